$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update the "Conversión del día" text block ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$new = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 10.3 = 42739.55 pesos`n✅ 42739.55 pesos = 10.28 = 976.73 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $new

# --- tasas sheet: update N10/O10/N12/O12 values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 97.09999999999999
$ws2.Range("O10").Value = 4150.01
$ws2.Range("N12").Value = 4156.99
$ws2.Range("O12").Value = 95
